# Replace the two "m:userdoc"/"m:enduserdoc" Word field codes that live in
# the document footer with their literal M2Doc text-tag equivalents, e.g.
#   <w:fldChar begin/><w:instrText>m:userdoc 'zone1'</w:instrText><w:fldChar end/>
# becomes
#   <w:t>{m:userdoc 'zone1'}</w:t>
# (same thing for "m:enduserdoc" -> "{m:enduserdoc}").
# Both fields sit alone in their own paragraph, so after deleting the field
# we re-insert the literal text into that now-empty paragraph.

$d = $word.ActiveDocument
$footer = $d.Sections.First.Footers.Item(1)

# --- Field 1: m:userdoc 'zone1'  (paragraph right after the intro text) ---
$fields = $footer.Range.Fields
$userdocField = $fields.Item(1)
$userdocField.Delete()

$p1 = $footer.Range.Duplicate
$p1.Start = 39
$p1.End = 39
$p1.InsertAfter("{m:userdoc 'zone1'}")

# --- Field 2: m:enduserdoc  (paragraph right after the table) ---
$fields2 = $footer.Range.Fields
$enddocField = $fields2.Item(1)
$enddocField.Delete()

$p2 = $footer.Range.Duplicate
$p2.Start = 79
$p2.End = 79
$p2.InsertAfter("{m:enduserdoc}")
